$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common (unchanged across every data row in this block)
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$prodId    = 100108
$producto  = "Tropicales y subtropicales"
$catId     = 100108005
$categoria = "Piña"
$variedad  = "Caramelo"
$origen    = "Ecuador"

function Fill-Row($r, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom, $unidad, $pkg, $kgUnidad) {
    $ws.Cells.Item($r, 1).Value  = 1
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $prodId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $catId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $pmin
    $ws.Cells.Item($r, 15).Value = $pmax
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $pkg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

# Step 1: insert two new rows at 233-234 for the new weekly report (Fecha = 44889)
$ws.Rows("233:234").Insert()

Fill-Row 233 44889 "Segunda" 250 31000 32000 31400 "$/caja 14 unidades" 2243 14
Fill-Row 234 44889 "Tercera" 170 31000 32000 31529 "$/caja 16 unidades" 1971 16

# Step 2: the 44691 block (now at rows 239-241) is missing an "Especial" entry; insert it
$ws.Rows("239:239").Insert()

Fill-Row 239 44691 "Especial" 200 23000 24000 23500 "$/caja 10 unidades" 2350 10
